$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; existing D:K data shifts to E:L
$ws.Columns("D").Insert()

# Copy the number-format/font/alignment from the (now-shifted) neighboring
# column E onto the freshly inserted column D cells so they match the
# rest of the table (date style for header rows, number style elsewhere).
$ws.Range("E7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E80").Copy()
$ws.Range("D80").PasteSpecial(-4122)

$ws.Range("E8:E35").Copy()
$ws.Range("D8:D35").PasteSpecial(-4122)
$ws.Range("E39:E77").Copy()
$ws.Range("D39:D77").PasteSpecial(-4122)
$ws.Range("E81:E102").Copy()
$ws.Range("D81:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the FY2018 (period ending 2018-12-31) data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 212400
$ws.Range("D9").Value = 147600
$ws.Range("D10").Value = 64800
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 108700
$ws.Range("D17").Value = 299800
$ws.Range("D18").Value = -87400
$ws.Range("D20").Value = 2200
$ws.Range("D21").Value = 23500
$ws.Range("D22").Value = 63600
$ws.Range("D23").Value = -148700
$ws.Range("D24").Value = -29600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -119100
$ws.Range("D27").Value = -119100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2200
$ws.Range("D33").Value = -119100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -119100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 224900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 54900
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 19800
$ws.Range("D46").Value = 299600
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 2434800
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 30200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2764600
$ws.Range("D57").Value = 26800
$ws.Range("D58").Value = 96300
$ws.Range("D59").Value = 38100
$ws.Range("D60").Value = 161200
$ws.Range("D61").Value = 1123600
$ws.Range("D62").Value = 171800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1456700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 549500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1307900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -119100
$ws.Range("D83").Value = 108700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -42400
$ws.Range("D91").Value = -52600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -52500
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 133800
$ws.Range("D101").Value = -800
$ws.Range("D102").Value = 38100
